# Auto-generated edit script: apply Asura_Profits.xlsx sheet updates
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 295.25
$ws.Range("I6").Value = 295.25
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 885.75
$ws.Range("L6").Value = 0
$ws.Range("M6").Value = -773.75
$ws.Range("N6").Value = $null

$ws.Range("H11").Value = 847
$ws.Range("I11").Value = 847
$ws.Range("K11").Value = 847
$ws.Range("M11").Value = -707

$ws.Range("H31").Value = 33270.43
$ws.Range("I31").Value = 33270.43
$ws.Range("K31").Value = 99811.29000000001
$ws.Range("M31").Value = -99581.29000000001

$ws.Range("H42").Value = 331.70587
$ws.Range("I42").Value = 223.28572
$ws.Range("J42").Value = 407.6
$ws.Range("K42").Value = 669.85716
$ws.Range("L42").Value = 1222.8
$ws.Range("M42").Value = -439.85716
$ws.Range("N42").Value = -1682.8

$ws.Range("H116").Value = 25002324
$ws.Range("I116").Value = 28573540
$ws.Range("J116").Value = 3800
$ws.Range("K116").Value = 28573540
$ws.Range("L116").Value = 3800
$ws.Range("M116").Value = -28570098
$ws.Range("N116").Value = -10684

$ws.Range("H138").Value = 2686.037
$ws.Range("I138").Value = 1671.9524
$ws.Range("J138").Value = 3331.3635
$ws.Range("K138").Value = 5015.857199999999
$ws.Range("L138").Value = 9994.0905
$ws.Range("M138").Value = 124.1428000000005
$ws.Range("N138").Value = -20274.0905

$ws.Range("H141").Value = 10557.714
$ws.Range("I141").Value = 3691.4285
$ws.Range("J141").Value = 17424
$ws.Range("K141").Value = 11074.2855
$ws.Range("L141").Value = 52272
$ws.Range("M141").Value = -5894.2855
$ws.Range("N141").Value = -62632

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 1144.9445
$ws.Range("I74").Value = 873.9
$ws.Range("K74").Value = 873.9
$ws.Range("M74").Value = 0.1000000000000227

$ws.Range("H77").Value = 1144.9445
$ws.Range("I77").Value = 873.9
$ws.Range("K77").Value = 4369.5
$ws.Range("M77").Value = -1.5

$ws.Range("H125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("L125").Value = 0
$ws.Range("N125").Value = $null

$ws.Range("H132").Value = 3473.356
$ws.Range("I132").Value = 3887.35
$ws.Range("K132").Value = 11662.05
$ws.Range("M132").Value = -9132.049999999999

$ws.Range("H134").Value = 33747.5
$ws.Range("J134").Value = 33747.5
$ws.Range("L134").Value = 33747.5
$ws.Range("N134").Value = -43887.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2300.9546
$ws.Range("I31").Value = 1826.6
$ws.Range("J31").Value = 3317.4285
$ws.Range("K31").Value = 1826.6
$ws.Range("L31").Value = 3317.4285
$ws.Range("M31").Value = -1531.6
$ws.Range("N31").Value = -3907.4285

$ws.Range("H34").Value = 2300.9546
$ws.Range("I34").Value = 1826.6
$ws.Range("J34").Value = 3317.4285
$ws.Range("K34").Value = 1826.6
$ws.Range("L34").Value = 3317.4285
$ws.Range("M34").Value = -1624.6
$ws.Range("N34").Value = -3721.4285

$ws.Range("H58").Value = 773341.6
$ws.Range("I58").Value = 1278762.4
$ws.Range("K58").Value = 1278762.4
$ws.Range("M58").Value = -1278559.4

$ws.Range("H68").Value = 32000
$ws.Range("J68").Value = 32000
$ws.Range("L68").Value = 32000
$ws.Range("N68").Value = -33498

$ws.Range("H71").Value = 32000
$ws.Range("J71").Value = 32000
$ws.Range("L71").Value = 96000
$ws.Range("N71").Value = -103488

$ws.Range("H86").Value = 4081.2
$ws.Range("I86").Value = 2937.3333
$ws.Range("J86").Value = 4571.4287
$ws.Range("K86").Value = 2937.3333
$ws.Range("L86").Value = 4571.4287
$ws.Range("M86").Value = -1814.3333
$ws.Range("N86").Value = -6817.4287

$ws.Range("H89").Value = 4081.2
$ws.Range("I89").Value = 2937.3333
$ws.Range("J89").Value = 4571.4287
$ws.Range("K89").Value = 14686.6665
$ws.Range("L89").Value = 22857.1435
$ws.Range("M89").Value = -9070.666499999999
$ws.Range("N89").Value = -34089.14350000001

$ws.Range("H107").Value = 348.63635
$ws.Range("I107").Value = 326.8
$ws.Range("J107").Value = 395.42856
$ws.Range("K107").Value = 326.8
$ws.Range("L107").Value = 395.42856
$ws.Range("M107").Value = 1593.2
$ws.Range("N107").Value = -4235.42856

$ws.Range("H136").Value = 773341.6
$ws.Range("I136").Value = 1278762.4
$ws.Range("K136").Value = 3836287.2
$ws.Range("M136").Value = -3833737.2

$ws.Range("H141").Value = 26431.6
$ws.Range("J141").Value = 27701.777
$ws.Range("L141").Value = 27701.777
$ws.Range("N141").Value = -38061.777

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H137").Value = 2296.9092
$ws.Range("I137").Value = 1345.7142
$ws.Range("J137").Value = 3961.5
$ws.Range("K137").Value = 4037.1426
$ws.Range("L137").Value = 11884.5
$ws.Range("M137").Value = 1062.8574
$ws.Range("N137").Value = -22084.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5899.75
$ws.Range("J70").Value = 6100
$ws.Range("L70").Value = 6100
$ws.Range("N70").Value = -6640

$ws.Range("H73").Value = 5899.75
$ws.Range("J73").Value = 6100
$ws.Range("L73").Value = 6100
$ws.Range("N73").Value = -7972

$ws.Range("H109").Value = 9223.0625
$ws.Range("J109").Value = 9223.0625
$ws.Range("L109").Value = 9223.0625
$ws.Range("N109").Value = -11303.0625

$ws.Range("H131").Value = 48769
$ws.Range("J131").Value = 48769
$ws.Range("L131").Value = 48769
$ws.Range("N131").Value = -58849

$ws.Range("H132").Value = 1735.234
$ws.Range("I132").Value = 1187.8572
$ws.Range("K132").Value = 3563.5716
$ws.Range("M132").Value = -1033.5716

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H36").Value = 500715
$ws.Range("J36").Value = 500715
$ws.Range("L36").Value = 500715
$ws.Range("N36").Value = -501839

$ws.Range("H131").Value = 0
$ws.Range("J131").Value = 0
$ws.Range("L131").Value = 0
$ws.Range("N131").Value = $null

$ws.Range("H132").Value = 6950.4707
$ws.Range("I132").Value = 8014.5454
$ws.Range("K132").Value = 24043.6362
$ws.Range("M132").Value = -21513.6362

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H123").Value = 24476
$ws.Range("J123").Value = 24476
$ws.Range("L123").Value = 24476
$ws.Range("N123").Value = -34276
